$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.436.43"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "1.571.45"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.55"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3694"
$ws.Range("E7").Value = "  +2.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.34"
$ws.Range("E8").Value = "  -2.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3317"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.162"
$ws.Range("E10").Value = "  +3.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07491"
$ws.Range("E11").Value = "  +1.39%  "

$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.74"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.933"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.892"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").Value = "1.557.79"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.82"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06720"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.421"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.45"
$ws.Range("E22").Value = "  +1.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.96"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("D24").Value = "22.418.87"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.371"
$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.625"
$ws.Range("E26").Value = "  +3.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.90"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.54"
$ws.Range("E28").Value = "  +1.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.941"
$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.67"
$ws.Range("E30").Value = "  +1.05%  "

$ws.Range("D31").Value = "1.736.16"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.088"
$ws.Range("E33").Value = "  -0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.975"
$ws.Range("E34").Value = "  -1.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.894"
$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08316"
$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("E37").Value = "  +1.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06386"
$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.298"
$ws.Range("E39").Value = "  +0.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2210"
$ws.Range("E40").Value = "  -0.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.322"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.37"
$ws.Range("E42").Value = "  +2.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6202"
$ws.Range("E43").Value = "  +2.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.03"
$ws.Range("E45").Value = "  +1.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6007"
$ws.Range("E46").Value = "  +3.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.768"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.040"
$ws.Range("E48").Value = "  +1.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.86"
$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.198"
$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07186"
$ws.Range("E51").Value = "  -0.29%  "
